# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450 — update
# the StructureDefinition "description" IG metadata sheet:
#   - Version 5.0.0 -> 6.0.0
#   - Date updated
#   - Publisher now has a value ("Alvearie Team")
#   - Duplicate "Contact" row replaced by a new "Jurisdiction" row
#   - Elements sheet: root Extension's Short/Definition updated to match
#     the new top-level Description

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)   # "Metadata" sheet

$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"

# Remove the second (duplicate) "Contact" / "No display for ContactDetail" row.
$ws.Rows.Item(11).Delete()

# The remaining "Contact" row (now row 10) becomes the new "Jurisdiction" row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

$ws2 = $wb.Worksheets.Item(2)  # "Elements" sheet

# Update the root Extension element's Short/Definition columns.
$ws2.Range("K2").Value = "Description"
$ws2.Range("L2").Value = "Plain text description that explains the insight score result"
